$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value = "Maule"
$ws.Cells.Item($row, 4).Value = 44890
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = "Arándano (blue)"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 3600
$ws.Cells.Item($row, 15).Value = 3600
$ws.Cells.Item($row, 16).Value = 3600
$ws.Cells.Item($row, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row, 19).Value = 1800
$ws.Cells.Item($row, 20).Value = 2
